$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force the Price/Volume columns to stay text-typed (they are stored as inline
# strings in the source, e.g. "1.00" / "70.258.89" / "  -0.88%  ") so assigning
# numeric-looking strings via COM does not silently coerce them to numbers.
$ws.Range("D2:E51").NumberFormat = "@"

$changes = @{
    "D2" = "70.258.89"
    "E2" = "  -0.88%  "
    "D3" = "3.485.46"
    "E3" = "  -2.04%  "
    "D4" = "1.00"
    "E4" = "  +0.04%  "
    "D5" = "616.58"
    "E5" = "  +2.47%  "
    "D6" = "168.74"
    "E6" = "  -1.85%  "
    "D7" = "3.481.37"
    "E7" = "  -2.01%  "
    "D8" = "0.601"
    "E8" = "  -2.19%  "
    "E9" = "  +0.14%  "
    "D10" = "0.195"
    "E10" = "  +0.56%  "
    "D11" = "7.15"
    "E11" = "  -2.85%  "
    "D12" = "0.571"
    "E12" = "  -2.66%  "
    "D13" = "45.04"
    "E13" = "  -2.77%  "
    "D14" = "0.0000270"
    "E14" = "  -2.27%  "
    "D15" = "4.055.89"
    "E15" = "  -1.85%  "
    "D16" = "8.27"
    "E16" = "  -0.71%  "
    "D17" = "594.97"
    "E17" = "  -2.53%  "
    "D18" = "70.394.70"
    "E18" = "  -0.68%  "
    "D19" = "3.494.08"
    "E19" = "  -1.78%  "
    "D20" = "0.121"
    "E20" = "  +1.51%  "
    "D21" = "17.37"
    "E21" = "  +0.22%  "
    "D22" = "0.866"
    "E22" = "  -1.43%  "
    "D23" = "8.87"
    "E23" = "  -4.67%  "
    "B24" = "Litecoin"
    "C24" = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
    "D24" = "96.36"
    "E24" = "  -0.34%  "
    "B25" = "InternetComputer(DFINITY)"
    "C25" = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
    "D25" = "15.28"
    "E25" = "  -2.74%  "
    "D26" = "3.65"
    "E26" = "  -1.60%  "
    "E27" = "  +0.04%  "
    "D28" = "2.50"
    "E28" = "  -3.85%  "
    "D29" = "33.31"
    "E29" = "  -2.00%  "
    "D30" = "8.80"
    "E30" = "  -2.86%  "
    "D31" = "7.97"
    "E31" = "  -3.11%  "
    "D32" = "2.86"
    "E32" = "  -6.57%  "
    "D33" = "1.26"
    "E33" = "  -2.66%  "
    "D34" = "6.64"
    "E34" = "  -6.02%  "
    "D35" = "576.59"
    "E35" = "  -18.32%  "
    "D36" = "0.0491"
    "E36" = "  +2.59%  "
    "D37" = "10.74"
    "E37" = "  -0.17%  "
    "D38" = "0.0974"
    "E38" = "  -3.13%  "
    "B39" = "FirstDigitalUSD"
    "C39" = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
    "D39" = "1.00"
    "E39" = "  +0.37%  "
    "B40" = "OKB"
    "C40" = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
    "D40" = "56.46"
    "E40" = "  -0.87%  "
    "D41" = "0.142"
    "E41" = "  -1.04%  "
    "D42" = "3.26"
    "E42" = "  -10.06%  "
    "D43" = "3.295.16"
    "E43" = "  -2.24%  "
    "D44" = "0.0₃0703"
    "E44" = "  +0.73%  "
    "D45" = "0.302"
    "E45" = "  -4.96%  "
    "D46" = "31.19"
    "E46" = "  -4.34%  "
    "D47" = "2.78"
    "E47" = "  -4.66%  "
    "D48" = "2.43"
    "E48" = "  -6.54%  "
    "E49" = "  -1.94%  "
    "D50" = "133.59"
    "E50" = "  -0.15%  "
    "E51" = "  -0.03%  "
}

foreach ($ref in $changes.Keys) {
    $ws.Range($ref).Value2 = $changes[$ref]
}

# Restore the default (unstyled) look now that the text values are committed.
$ws.Range("D2:E51").Style = "Normal"
